$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.791.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.723.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '657.15'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.427'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("E9").Value = '  +3.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.722.32'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.413.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.52%  '
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.676.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +17.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.718.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.533'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '524.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000205'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '103.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("E29").Value = '  -3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.13%  '
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("E33").Value = '  +13.16%  '
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.30%  '
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '658.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.608'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.19%  '
$ws.Range("E39").Value = '  +3.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +19.27%  '
$ws.Range("E41").Value = '  +5.05%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.982'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.18%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +22.80%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.36%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.451'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0460'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.66'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.51%  '
